# update file đánh giá
# - Append the week-6 "Hoàn thành chức năng hiển thị sản phẩm và thông tin"
#   work-log entry (5h, MSSV 1359037 / Trần Long Sơn) to both the
#   "Tổng quan" summary sheet (row 21) and the "Chi tiết" detail sheet
#   (row 24, dated 10/6/2016).
# - Leave the workbook with the "Tổng quan" sheet active/selected again.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Tổng quan"
$ws2 = $wb.Worksheets.Item(2)   # "Chi tiết"

# --- Sheet "Tổng quan": new row 21 -----------------------------------
$ws1.Range("A21").Value = 20
$ws1.Range("B21").Value = 1359037
$ws1.Range("C21").Value = "Trần Long Sơn "
$ws1.Range("D21").Value = 5
$ws1.Range("E21").Value = "Hoàn thành chức năng hiển thị sản phẩm và thông tin"

# --- Sheet "Chi tiết": new row 24 -------------------------------------
$ws2.Range("A24").Value = "project "

# Match the date formatting already used by the row above it (B23)
# instead of inventing a brand new number format.
$ws2.Range("B23").Copy()
$ws2.Range("B24").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("B24").Value = 42649          # 10/6/2016

$ws2.Range("C24").Value = "Hoàn thành chức năng hiển thị sản phẩm và thông tin"
$ws2.Range("D24").Value = 1359037
$ws2.Range("E24").Value = 5

# --- Selection / active sheet -----------------------------------------
# Select "Chi tiết" first (so its selection becomes C24, the new row),
# then select "Tổng quan" last so it ends up the active/visible tab.
$ws2.Range("C24").Select()
$ws1.Range("E23").Select()
